# Append 8 new NBA 2023-24 game rows (rows 664-671) to Sheet1, continuing
# the existing table of games (Away team, Away Pts, Home team, Home Pts,
# Overtime, Attend., Arena, Win, Loss).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$games = @(
    @{ Away = "Dallas Mavericks";       AwayPts = 148; Home = "Atlanta Hawks";         HomePts = 143; OT = "No"; Attend = 17832; Arena = "State Farm Arena";         Win = "Dallas Mavericks";       Loss = "Atlanta Hawks" },
    @{ Away = "Houston Rockets";        AwayPts = 138; Home = "Charlotte Hornets";     HomePts = 104; OT = "No"; Attend = 17832; Arena = "Spectrum Center";          Win = "Houston Rockets";        Loss = "Charlotte Hornets" },
    @{ Away = "Phoenix Suns";           AwayPts = 131; Home = "Indiana Pacers";        HomePts = 133; OT = "No"; Attend = 17832; Arena = "Gainbridge Fieldhouse";    Win = "Indiana Pacers";         Loss = "Phoenix Suns" },
    @{ Away = "Los Angeles Clippers";   AwayPts = 127; Home = "Toronto Raptors";       HomePts = 107; OT = "No"; Attend = 17832; Arena = "Scotiabank Arena";         Win = "Los Angeles Clippers";   Loss = "Toronto Raptors" },
    @{ Away = "Orlando Magic";          AwayPts = 106; Home = "Memphis Grizzlies";     HomePts = 107; OT = "No"; Attend = 17832; Arena = "FedEx Forum";              Win = "Memphis Grizzlies";      Loss = "Orlando Magic" },
    @{ Away = "Cleveland Cavaliers";    AwayPts = 112; Home = "Milwaukee Bucks";       HomePts = 100; OT = "No"; Attend = 17832; Arena = "Fiserv Forum";             Win = "Cleveland Cavaliers";    Loss = "Milwaukee Bucks" },
    @{ Away = "Oklahoma City Thunder";  AwayPts = 107; Home = "New Orleans Pelicans";  HomePts = 83;  OT = "No"; Attend = 17832; Arena = "Smoothie King Center";     Win = "Oklahoma City Thunder";  Loss = "New Orleans Pelicans" },
    @{ Away = "Portland Trail Blazers"; AwayPts = 100; Home = "San Antonio Spurs";     HomePts = 116; OT = "No"; Attend = 17832; Arena = "Frost Bank Center";        Win = "San Antonio Spurs";      Loss = "Portland Trail Blazers" }
)

$startRow = 664
for ($i = 0; $i -lt $games.Count; $i++) {
    $row = $startRow + $i
    $g = $games[$i]

    $ws.Cells.Item($row, 1).Value = $g.Away
    $ws.Cells.Item($row, 2).Value = $g.AwayPts
    $ws.Cells.Item($row, 2).NumberFormat = "#,##0"
    $ws.Cells.Item($row, 3).Value = $g.Home
    $ws.Cells.Item($row, 4).Value = $g.HomePts
    $ws.Cells.Item($row, 4).NumberFormat = "#,##0"
    $ws.Cells.Item($row, 5).Value = $g.OT
    $ws.Cells.Item($row, 6).Value = $g.Attend
    $ws.Cells.Item($row, 7).Value = $g.Arena
    $ws.Cells.Item($row, 8).Value = $g.Win
    $ws.Cells.Item($row, 9).Value = $g.Loss
}
